# Update column O (player_id_y) values on Sheet1 to reflect the
# refreshed player_id_y values from the updated Kaggle CSV export.
# Mapping is row -> new value (only rows that actually changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = 3838
    3  = 3418
    4  = 5020
    5  = 2870
    6  = 3782
    7  = 3360
    8  = 3360
    9  = 1681
    10 = 2336
    11 = 3360
    12 = 2336
    13 = 2224
    14 = 2358
    15 = 3198
    16 = 2336
    17 = 2431
    18 = 3441
    19 = 3269
    20 = 3715
    24 = 1024
    26 = 4346
    27 = 1175
    28 = 1146
    29 = 2799
    30 = 4908
    32 = 1291
    33 = 979
    34 = 1357
    35 = 1357
    36 = 4295
    37 = 1624
    38 = 4432
    39 = 4295
    41 = 3112
    42 = 3112
}

foreach ($row in $updates.Keys) {
    $ws.Range("O$row").Value = $updates[$row]
}
